$d = $word.ActiveDocument

function Insert-ItalicParagraphAfter($matchText, $newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($matchText + "`r")) {
            $p.Range.InsertParagraphAfter()
            $newPara = $d.Paragraphs.Item($i + 1)
            $newRange = $newPara.Range
            $newRange.InsertAfter($newText)
            $newRange2 = $newPara.Range
            $newRange2.MoveEnd(1, -1)
            $newRange2.Font.Italic = $true
            return
        }
    }
}

# 1. Update the activation date
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2. Insert English translation after the "Objetivos" paragraph (Portuguese)
Insert-ItalicParagraphAfter `
    "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte." `
    "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."

# 3. Insert English translation after the "Programa resumido" paragraph (Portuguese)
Insert-ItalicParagraphAfter `
    "A definir, de acordo com o tópico programado." `
    "To be defined, according to the programmed topic."

# 4. Insert English translation after the "Programa" paragraph (Portuguese)
Insert-ItalicParagraphAfter `
    "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação." `
    "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."

Write-Output "done"
